$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 222, pushing existing data (rows 222-261) down to rows 224-263.
$ws.Rows.Item(222).Resize(2).Insert()

# New row 222: Patagonia, 1a (guarda)
$ws.Range("A222").Value = 4
$ws.Range("B222").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C222").Value = "Los Lagos"
$ws.Range("D222").Value = 44504
$ws.Range("D222").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E222").Value = 10
$ws.Range("F222").Value = 100114001
$ws.Range("G222").Value = "Papa"
$ws.Range("H222").Value = "Patagonia"
$ws.Range("I222").Value = "1a (guarda)"
$ws.Range("J222").Value = 150
$ws.Range("K222").Value = 9000
$ws.Range("L222").Value = 9000
$ws.Range("M222").Value = 9000
$ws.Range("N222").Value = '$/saco 25 kilos'
$ws.Range("O222").Value = "Provincia de Llanquihue"
$ws.Range("P222").Value = 360
$ws.Range("Q222").Value = 25
$ws.Range("R222").Value = "Hortaliza"

# New row 223: Pehuenche, 1a nueva(o)
$ws.Range("A223").Value = 4
$ws.Range("B223").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C223").Value = "Los Lagos"
$ws.Range("D223").Value = 44504
$ws.Range("D223").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E223").Value = 10
$ws.Range("F223").Value = 100114001
$ws.Range("G223").Value = "Papa"
$ws.Range("H223").Value = "Pehuenche"
$ws.Range("I223").Value = "1a nueva(o)"
$ws.Range("J223").Value = 150
$ws.Range("K223").Value = 18000
$ws.Range("L223").Value = 18000
$ws.Range("M223").Value = 18000
$ws.Range("N223").Value = '$/saco 25 kilos'
$ws.Range("O223").Value = "Región de La Araucanía"
$ws.Range("P223").Value = 720
$ws.Range("Q223").Value = 25
$ws.Range("R223").Value = "Hortaliza"
